# Update the generated non-convex experiment values (MorganPatrone2006a,
# "Strong Stationary" generator, alpha = 0): a new run of the generator
# produced a new random point (x, y) and refreshed every value that is
# derived from it across the workbook's sheets.
#
# All of these "numeric" values are stored in the workbook as plain text
# (shared strings), not as real numbers, so we must write them in a way
# that keeps Excel from auto-converting them into numeric cells (which
# would also pull in an unwanted number-format style). Writing the value
# through a temporary formula and then doing a Copy + PasteSpecial
# (values only) reliably yields a genuine text cell with no style churn.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Sheet: Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $ws2.Range("A2") "2.8499999999999996 - x"
Set-TextValue $ws2.Range("B2") "-3.3499999999999996"
Set-TextValue $ws2.Range("D2") "0.3"
Set-TextValue $ws2.Range("A3") "-2.8499999999999996 + x"
Set-TextValue $ws2.Range("B3") "2.3499999999999996"
Set-TextValue $ws2.Range("D3") "0.09"

# --- Sheet: Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws3.Range("A2") "4.449999999999999 - y"
Set-TextValue $ws3.Range("B2") "-5.449999999999999"
Set-TextValue $ws3.Range("D2") "0.19"
Set-TextValue $ws3.Range("E2") "4.0"
Set-TextValue $ws3.Range("F2") "0.1"
Set-TextValue $ws3.Range("A3") "-4.449999999999999 + y"
Set-TextValue $ws3.Range("B3") "3.4499999999999993"
Set-TextValue $ws3.Range("D3") "0.79"
Set-TextValue $ws3.Range("E3") "2.9"
Set-TextValue $ws3.Range("F3") "7.800000000000001"

# --- Sheet: Punto_modificado ---
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "2.8499999999999996"
Set-TextValue $ws4.Range("B2") "4.449999999999999"

# --- Sheet: Vector_bf ---
$ws5 = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws5.Range("A2") "-3.4499999999999997"

# --- Sheet: Vector_BF ---
$ws6 = $wb.Worksheets.Item("Vector_BF")
Set-TextValue $ws6.Range("A2") "1.21"
Set-TextValue $ws6.Range("A3") "2.1"
